# Applies the target edit to the presentation:
#   1. Three tables (on slides 14, 15, 16) get their table style switched
#      from {29D4EDF6-4927-4204-A4FB-1F661AD6060E} ("Table_0") to the
#      built-in style {47B1AC74-648B-4459-8E17-9B8E1E5C58C0}.
#   2. The presentation's theme color scheme (currently the "Integral" /
#      "Red Violet" palette) is switched back to the standard "Office"
#      palette.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{47B1AC74-648B-4459-8E17-9B8E1E5C58C0}"

14, 15, 16 | ForEach-Object {
    $slide = $p.Slides.Item($_)
    $tbl = $slide.Shapes.Item(1).Table
    $tbl.ApplyStyle($newTableStyle)
}

# --- 2. Theme colors ---------------------------------------------------
# Restore the default Office color scheme (dk1,lt1,dk2,lt2,accent1-6,
# hlink,folHlink) on the slide master's color scheme.
function HexToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexToBgr $officeColors[$i]
}
